$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C9 and C10 values, and give C9 the same style as C10 (s="1")
$ws.Range("C10").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = 3.53
$ws.Range("C10").Value = 4.52

# Populate rows 11-50 with new grocery items
$ws.Range("A11").Value = "Eggs"
$ws.Range("B11").Value = 15246853.0
$ws.Range("C11").Value = 3.48
$ws.Range("D11").Value = 4.0
$ws.Range("A12").Value = "Granola"
$ws.Range("B12").Value = 15378954.0
$ws.Range("C12").Value = 2.29
$ws.Range("D12").Value = 1.0
$ws.Range("A13").Value = "Brussel Sprouts"
$ws.Range("B13").Value = 152322976.0
$ws.Range("C13").Value = 5.11
$ws.Range("D13").Value = 2.0
$ws.Range("A14").Value = "Steak"
$ws.Range("B14").Value = 1424.0
$ws.Range("C14").Value = 5.89
$ws.Range("D14").Value = 3.0
$ws.Range("A15").Value = "Sausage"
$ws.Range("B15").Value = 1643397.0
$ws.Range("C15").Value = 3.99
$ws.Range("D15").Value = 3.0
$ws.Range("A16").Value = "Cinnamon"
$ws.Range("B16").Value = 168905.0
$ws.Range("C16").Value = 2.26
$ws.Range("D16").Value = 5.0
$ws.Range("A17").Value = "Nutmeg"
$ws.Range("B17").Value = 954378.0
$ws.Range("C17").Value = 2.77
$ws.Range("D17").Value = 5.0
$ws.Range("A18").Value = "Oregano"
$ws.Range("B18").Value = 865588.0
$ws.Range("C18").Value = 1.34
$ws.Range("D18").Value = 5.0
$ws.Range("A19").Value = "Vanilla"
$ws.Range("B19").Value = 6784588.0
$ws.Range("C19").Value = 4.89
$ws.Range("D19").Value = 5.0
$ws.Range("A20").Value = "Quinoa"
$ws.Range("B20").Value = 486358.0
$ws.Range("C20").Value = 4.01
$ws.Range("D20").Value = 2.0
$ws.Range("A21").Value = "Peppers"
$ws.Range("B21").Value = 6886355.0
$ws.Range("C21").Value = 3.33
$ws.Range("D21").Value = 2.0
$ws.Range("A22").Value = "Toilet Paper"
$ws.Range("B22").Value = 3658964.0
$ws.Range("C22").Value = 10.22
$ws.Range("D22").Value = 6.0
$ws.Range("A23").Value = "Paper Towels"
$ws.Range("B23").Value = 3944588.0
$ws.Range("C23").Value = 18.74
$ws.Range("D23").Value = 6.0
$ws.Range("A24").Value = "Detergent"
$ws.Range("B24").Value = 555666577.0
$ws.Range("C24").Value = 9.69
$ws.Range("D24").Value = 6.0
$ws.Range("A25").Value = "Kielbasa"
$ws.Range("B25").Value = 8627838.0
$ws.Range("C25").Value = 4.23
$ws.Range("D25").Value = 3.0
$ws.Range("A26").Value = "Lettuce"
$ws.Range("B26").Value = 384589.0
$ws.Range("C26").Value = 4.74
$ws.Range("D26").Value = 2.0
$ws.Range("A27").Value = "Arugula"
$ws.Range("B27").Value = 644777.0
$ws.Range("C27").Value = 6.74
$ws.Range("D27").Value = 2.0
$ws.Range("A28").Value = "Chocolate"
$ws.Range("B28").Value = 6848864988.0
$ws.Range("C28").Value = 3.99
$ws.Range("D28").Value = 1.0
$ws.Range("A29").Value = "Ice Cream"
$ws.Range("B29").Value = 3688965.0
$ws.Range("C29").Value = 4.56
$ws.Range("D29").Value = 3.0
$ws.Range("A30").Value = "Tomatos"
$ws.Range("B30").Value = 684955.0
$ws.Range("C30").Value = 3.45
$ws.Range("D30").Value = 2.0
$ws.Range("A31").Value = "Parsnips"
$ws.Range("B31").Value = 697742.0
$ws.Range("C31").Value = 2.22
$ws.Range("D31").Value = 2.0
$ws.Range("A32").Value = "Turnips"
$ws.Range("B32").Value = 65786.0
$ws.Range("C32").Value = 3.22
$ws.Range("D32").Value = 2.0
$ws.Range("A33").Value = "Chives"
$ws.Range("B33").Value = 546987.0
$ws.Range("C33").Value = 4.11
$ws.Range("D33").Value = 2.0
$ws.Range("A34").Value = "Mushrooms"
$ws.Range("B34").Value = 68744.0
$ws.Range("C34").Value = 8.32
$ws.Range("D34").Value = 2.0
$ws.Range("A35").Value = "Carrots"
$ws.Range("B35").Value = 3888965.0
$ws.Range("C35").Value = 3.61
$ws.Range("D35").Value = 2.0
$ws.Range("A36").Value = "Juice"
$ws.Range("B36").Value = 56658.0
$ws.Range("C36").Value = 5.22
$ws.Range("D36").Value = 3.0
$ws.Range("A37").Value = "Tortillas"
$ws.Range("B37").Value = 6878555.0
$ws.Range("C37").Value = 4.96
$ws.Range("D37").Value = 1.0
$ws.Range("A38").Value = "Bread"
$ws.Range("B38").Value = 56987.0
$ws.Range("C38").Value = 2.06
$ws.Range("D38").Value = 1.0
$ws.Range("A39").Value = "Rolls"
$ws.Range("B39").Value = 999999.0
$ws.Range("C39").Value = 1.79
$ws.Range("D39").Value = 1.0
$ws.Range("A40").Value = "Donuts"
$ws.Range("B40").Value = 3654.0
$ws.Range("C40").Value = 6.21
$ws.Range("D40").Value = 1.0
$ws.Range("A41").Value = "Cake"
$ws.Range("B41").Value = 77728698.0
$ws.Range("C41").Value = 8.02
$ws.Range("D41").Value = 1.0
$ws.Range("A42").Value = "Dank Lobster"
$ws.Range("B42").Value = 65759.0
$ws.Range("C42").Value = 101.01
$ws.Range("D42").Value = 3.0
$ws.Range("A43").Value = "Dried Fruit"
$ws.Range("B43").Value = 668789.0
$ws.Range("C43").Value = 5.02
$ws.Range("D43").Value = 5.0
$ws.Range("A44").Value = "Gatorade"
$ws.Range("B44").Value = 367555.0
$ws.Range("C44").Value = 4.21
$ws.Range("D44").Value = 1.0
$ws.Range("A45").Value = "Mango Chutney"
$ws.Range("B45").Value = 333331.0
$ws.Range("C45").Value = 3.86
$ws.Range("D45").Value = 5.0
$ws.Range("A46").Value = "Nutella"
$ws.Range("B46").Value = 745886.0
$ws.Range("C46").Value = 4.63
$ws.Range("D46").Value = 5.0
$ws.Range("A47").Value = "Peanut Butter"
$ws.Range("B47").Value = 836848.0
$ws.Range("C47").Value = 3.22
$ws.Range("D47").Value = 5.0
$ws.Range("A48").Value = "Jelly"
$ws.Range("B48").Value = 6815896.0
$ws.Range("C48").Value = 3.33
$ws.Range("D48").Value = 5.0
$ws.Range("A49").Value = "Pickles"
$ws.Range("B49").Value = 465444.0
$ws.Range("C49").Value = 4.08
$ws.Range("D49").Value = 5.0
$ws.Range("A50").Value = "Kombucha"
$ws.Range("B50").Value = 11111111119.0
$ws.Range("C50").Value = 4.19
$ws.Range("D50").Value = 3.0

# Apply the same style (s="1") used by row 10 to the newly added rows 11-50
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D50").PasteSpecial(-4122)

